$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new match row (row 3) with the same shape/columns as row 2.
# Columns G:K hold numeric-looking values that must stay text (matching
# the existing t="str" cells in row 2), so force a Text number format
# before writing them - otherwise Excel would coerce "0"/"2"/"0.00" into
# real numbers.
$ws.Range("G3:K3").NumberFormat = "@"

$ws.Range("A3").Value = " Dubai (DSC)"
$ws.Range("B3").Value = " September 24 2020"
$ws.Range("C3").Value = "Kings XI won by 97 runs"
$ws.Range("D3").Value = "Royal Challengers Bangalore"
$ws.Range("E3").Value = "Kings XI Punjab"
$ws.Range("F3").Value = "Umesh Yadav "
$ws.Range("G3").Value = "0"
$ws.Range("H3").Value = "2"
$ws.Range("I3").Value = "0"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "0.00"
